$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert three new rows at the bottom of the table (rows 27-29) and fill
#    in the three new "match" related functions.
# ---------------------------------------------------------------------------
$ws.Rows("27:29").Insert()

# Row 27 - Get Free Time
$ws.Range("A27").Value = 25
$ws.Range("C27").Value = "Get Free Time"
$ws.Range("D27").Value = "/swp49x-ffrs/match/free-time"
$ws.Range("E27").Value = "GET"
$ws.Range("G27").Value = "field-owner-id, field-type-id, date"
$ws.Range("H27").Value = "List<TimeSlotEntity>"
$ws.Range("I27").Value = "200, OK"

# Row 28 - Reserve Friendly Match
$ws.Range("A28").Value = 26
$ws.Range("C28").Value = "Reserve Friendly Match"
$ws.Range("D28").Value = "/swp49x-ffrs/match/friendly-match"
$ws.Range("E28").Value = "POST"
$ws.Range("F28").Value = "InputReservationDTO"
$ws.Range("H28").Value = "FriendlyMatchEntity"
$ws.Range("I28").Value = "201, CREATED"

# Row 29 - Get Match Up Coming
$ws.Range("A29").Value = 27
$ws.Range("C29").Value = "Get Match Up Coming"
$ws.Range("D29").Value = "/swp49x-ffrs/match/upcoming-match"
$ws.Range("E29").Value = "GET"
$ws.Range("G29").Value = "field-owner-id, field-type-id, date"
$ws.Range("H29").Value = "List<TimeSlotEntity>"
$ws.Range("I29").Value = "200, OK"

# ---------------------------------------------------------------------------
# 2. Give the new rows the same banded fill/border used by the "orange"
#    group (rows 11-13) by copying the formatting of an existing row that
#    already carries it - this re-uses the existing style instead of
#    minting new duplicate style records.
# ---------------------------------------------------------------------------
$ws.Range("A11:I11").Copy()
$ws.Range("A27:I27").PasteSpecial(-4122)
$ws.Range("A28:I28").PasteSpecial(-4122)
$ws.Range("A29:I29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Row 26 previously had an inconsistent style (a leftover "no top
#    border" variant on some of its cells). Re-level the whole row onto the
#    uniform style used by the rest of the group (e.g. row 25).
# ---------------------------------------------------------------------------
$ws.Range("A25:I25").Copy()
$ws.Range("A26:I26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Row 21 (last row of the "pink" group) drops its bottom border so it no
#    longer reads as a fully boxed-in row.
# ---------------------------------------------------------------------------
$ws.Range("A21:I21").Borders(9).LineStyle = 0

# ---------------------------------------------------------------------------
# 5. Restore the selection to where the author left off editing.
# ---------------------------------------------------------------------------
$ws.Range("D18").Select()
